# PyramidPanic.xlsx edit:
#  - Duplicate "week 50" sheet, rename the copy to "week 2", place it right
#    after "week 50" (before "Totaal").
#  - Update the new sheet's log entry (row 7) with the new date/time/text,
#    and clear out the remaining (unused) log rows 8-17.
#  - Insert a new row into "Totaal" for week "2" referencing the new sheet,
#    pushing the old "Totaal:" summary row down one.
#  - Leave the new "week 2" sheet as the active/selected sheet, matching
#    the workbook's bookView state after the edit.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "week 50" into a new "week 2" sheet -----------------
$ws50 = $wb.Worksheets.Item("week 50")
$ws50.Copy($null, $ws50)
$weekTwo = $wb.Worksheets.Item("week 50 (2)")
$weekTwo.Name = "week 2"

# "week 50" is no longer the active tab; its selection also moved on to
# cover the whole used range (A1:G18) with G18 as the last-clicked cell.
$ws50.Activate()
$ws50.Range("A1:G18").Select()

# --- 2. Update the log entry in row 7 ----------------------------------
$weekTwo.Range("B7").Value = 41648
$weekTwo.Range("C7").Value = 0.36458333333333331
$weekTwo.Range("D7").Value = 0.38194444444444442
$weekTwo.Range("F7").Value = "Menu Gemaakt, hij ging uit het beeld. scoreScene en LoadScene aangemaakt en van elke scene kan je terug via b"

# --- 3. Clear out the rest of the (unused) log rows --------------------
$weekTwo.Range("C8:D17").ClearContents()
$weekTwo.Range("F8:F14").ClearContents()

# Row 7 now wraps a long activity description (auto row height); the rows
# that lost their text go back to the default height.
$weekTwo.Rows.Item(7).RowHeight = 42.75
$weekTwo.Rows.Item(8).EntireRow.AutoFit()
$weekTwo.Rows.Item(10).EntireRow.AutoFit()
$weekTwo.Rows.Item(11).EntireRow.AutoFit()
$weekTwo.Rows.Item(14).EntireRow.AutoFit()

# --- 4. Selection / active state for "week 2" ---------------------------
$weekTwo.Activate()
$weekTwo.Range("F7").Select()

# --- 5. Update "Totaal" sheet: insert a row for week "2" ---------------
$totaal = $wb.Worksheets.Item("Totaal")
$totaal.Rows.Item(10).Insert()
$totaal.Range("A10").Value = 2
$totaal.Range("B10").Formula = "='week 2'!G18"
$totaal.Range("B28").Select()

# --- 6. Leave "week 2" as the active sheet/tab --------------------------
$weekTwo.Activate()
$weekTwo.Range("F7").Select()
